$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A37").Value = "Serveur"
$ws.Range("B37").Value = "Client"
$ws.Range("C37").Value = "BE"
$ws.Range("D37").Value = "idWinnerTeam"
$ws.Range("E37").Value = "Le serveur indique qu'une équipe a gagné."

$ws.Range("E37").Select()
